# Update "想去人数" (F column) values on the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 117
$wsExhibition.Range("F4").Value = 1634
$wsExhibition.Range("F5").Value = 0
$wsExhibition.Range("F6").Value = 23
$wsExhibition.Range("F8").Value = 145
$wsExhibition.Range("F10").Value = 501

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 393
$wsAll.Range("F3").Value = 117
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 417
$wsAll.Range("F8").Value = 0
$wsAll.Range("F9").Value = 63
$wsAll.Range("F10").Value = 501
